$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "DTOS y Autommer" (AutoMapper) section starting at row 36
$ws.Range("B36").Value = "DTOS y Autommer"
$ws.Range("B37").Value = "Instalar autommaper.dependienci extension"
$ws.Range("B38").Value = "Crear carpeta utilidades para guardar la clase de automapper"
$ws.Range("B39").Value = "Instanciar automapper en el startup"
$ws.Range("B40").Value = "Crear la clase de autommaper "
$ws.Range("D40").Value = "   CreateMap<AutorCreacionDTO, Autor>();"

# New note next to "Validaciones por defecto" block (row 27)
$ws.Range("K27").Value = "ctrl + r + r "
$ws.Range("L27").Value = "Renombra variables"

# Move the active selection to where the new note was entered
$ws.Range("L27").Select()
